# The author removed one post entry ("「あなたは冷たい」...") that used to be
# row 81 in the data table. All the rows below it (formerly 82..179) shift up
# by one (81..178), and the sheet dimension shrinks from A1:C179 to A1:C178.
#
# Deleting the entire row 81 reproduces exactly that: the old row 82 content
# becomes the new row 81, etc., down to the old row 179 becoming the new row
# 178, with no row 179 left over.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(81).Delete()
